$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Calculation"
$ws.Range("C2").Value = "Hypotenuse"
$ws.Range("F2").Value = "Opposite"
$ws.Range("E2").Value = "Adjacent"
$ws.Range("G2").Value = "Angle 1"
$ws.Range("H2").Value = "Angle 2"

$ws.Range("A3").Value = 1
$ws.Range("C3").Value = 5
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 36.87
$ws.Range("H3").Value = 53.13

$ws.Range("A4").Value = 2
$ws.Range("C4").Value = "5mm"
$ws.Range("G4").Value = 47
$ws.Range("H4").Value = 43

$ws.Range("A5").Value = 3
$ws.Range("C5").Value = "11.4km"
$ws.Range("F5").Value = "9km"
$ws.Range("E5").Value = "7km"
$ws.Range("G5").Value = 52.13
$ws.Range("H5").Value = 37.87

$ws.Range("E4").Value = "3.41mm"
$ws.Range("F4").Value = "3.66mm"

$ws.Range("D2:D5").Value = $null

$ws.Range("F7").Select()
